$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("summary_result")
$ws2 = $wb.Worksheets.Item("summary")

# --- Core data edit: TC14 ReviewAcademicServices / TC15 EditReview test-data counts
#     increased from 19 to 28 (Pass stays 0, Fail tracks Test Data) on both the
#     "summary_result" sheet (rows 23-24) and the mirrored "summary" sheet (rows 15-16).
$ws1.Range("B23").Value = 28
$ws1.Range("D23").Value = 28
$ws1.Range("B24").Value = 28
$ws1.Range("D24").Value = 28

$ws2.Range("B15").Value = 28
$ws2.Range("D15").Value = 28
$ws2.Range("B16").Value = 28
$ws2.Range("D16").Value = 28

# Rows 15-16 on "summary" were reformatted to match the lighter (non-bold) row
# style used elsewhere on that sheet (e.g. row 3) instead of the bold/boxed style.
$ws2.Range("B3:D3").Copy()
$ws2.Range("B15:D16").PasteSpecial(-4122)
$ws2.Range("E3:F3").Copy()
$ws2.Range("E15:F16").PasteSpecial(-4122)

# --- Number-format touch-ups (2 decimal places) on several %Pass/%Fail cells
$ws1.Range("E41:F41").NumberFormat = "0.00"
$ws1.Range("E46:F46").NumberFormat = "0.00"
$ws1.Range("E51:F51").NumberFormat = "0.00"
$ws2.Range("E31:F31").NumberFormat = "0.00"

# --- View state: user ended up with "summary" as the active sheet/tab, having
#     scrolled/selected E41:F41 on "summary_result" and E31:F31 on "summary".
$ws1.Activate()
$ws1.Range("E41:F41").Select()

$ws2.Activate()
$ws2.Range("E31:F31").Select()
